$d = $word.ActiveDocument

$replacements = @(
    @("24×82=1968", "72×69=4968"),
    @("27×11=297",  "65×42=2730"),
    @("98×65=6370", "70×44=3080"),
    @("65×91=5915", "68×66=4488"),
    @("94×25=2350", "61×50=3050"),
    @("75×11=825",  "11×57=627"),
    @("23×28=644",  "78×30=2340"),
    @("38×69=2622", "51×18=918"),
    @("43×77=3311", "35×82=2870"),
    @("92×68=6256", "76×37=2812"),
    @("94×18=1692", "18×86=1548"),
    @("39×49=1911", "97×21=2037"),
    @("87×46=4002", "52×13=676"),
    @("36×54=1944", "63×19=1197"),
    @("58×97=5626", "34×43=1462"),
    @("47×50=2350", "94×72=6768"),
    @("59×39=2301", "67×57=3819"),
    @("96×48=4608", "56×58=3248"),
    @("68×74=5032", "29×80=2320"),
    @("43×17=731",  "28×32=896"),
    @("11×78=858",  "50×37=1850"),
    @("25×45=1125", "88×50=4400"),
    @("40×75=3000", "40×62=2480"),
    @("66×78=5148", "75×70=5250"),
    @("61×28=1708", "80×69=5520")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
